# DOSINZAGE2-445: Updated test data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill the new "Voorvoegsels" column (D) with the empty-value placeholder
# for every data row that doesn't have a real name prefix.
foreach ($row in 2,3,4,5,6,8) {
    $ws.Cells.Item($row, 4).Value = "[Leeg]"
}

# Row 8's "Onderzoek.Verrichting.VerrichtingType" cell (I8) was blank;
# it now also gets the placeholder value.
$ws.Range("I8").Value = "[Leeg]"

# Update the active selection to reflect where the editor left off.
$ws.Activate()
$ws.Range("D16").Select()
